$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "as of" quarter header in BB1 (copy date style from BA1, then set the new date)
$ws.Range("BA1").Copy($ws.Range("BB1"))
$ws.Range("BB1").Value = 45986

# Most vintages (rows 2-81) simply repeat their last known (BA) value in the new BB column
for ($r = 2; $r -le 81; $r++) {
    $ws.Range("BA$r").Copy($ws.Range("BB$r"))
}

# The two most recent existing vintages get revised values in the new column
$ws.Range("BB82").Value = -0.1118837721692358
$ws.Range("BB83").Value = 0.3266766184601977

# Brand new vintage row for this quarter
$ws.Range("A83").Copy($ws.Range("A84"))
$ws.Range("A84").Value = 45884
$ws.Range("BB84").Value = 0.325608361860148
